$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.450.49"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.374.63"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.66"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.46"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.84"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.50"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.738.33"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.35"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.393.63"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.529.23"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.94"
$ws.Range("E19").Value = "  +6.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  -5.09%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.45"
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.86"
$ws.Range("E24").Value = "  -3.73%  "
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.17"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -5.53%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0987"
$ws.Range("E30").Value = "  +5.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.37"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.26"
$ws.Range("E32").Value = "  -6.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "166.98"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("E38").Value = "  +10.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.01"
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.00"
$ws.Range("E42").Value = "  -6.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.97"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  -5.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.74"
$ws.Range("E46").Value = "  -8.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.820.27"
$ws.Range("E47").Value = "  +9.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.22"
$ws.Range("E48").Value = "  +5.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.81"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.14"
$ws.Range("E51").Value = "  -6.81%  "
